$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D's text formatting (avoid Excel auto-converting
# numeric-looking strings like "22.12" or "0.0500" into floating point
# numbers, which would lose exact text representation / trailing zeros).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '62.003.72'
$ws.Range("D3").Value = '3.416.70'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '409.01'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = '129.17'
$ws.Range("E6").Value = '  -4.37%  '
$ws.Range("D7").Value = '0.638'
$ws.Range("E7").Value = '  +7.99%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.734'
$ws.Range("E9").Value = '  +7.08%  '
$ws.Range("E10").Value = '  +17.03%  '
$ws.Range("D11").Value = '42.45'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '0.0000217'
$ws.Range("E12").Value = '  +67.15%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").Value = '3.960.94'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = '8.99'
$ws.Range("E15").Value = '  +6.58%  '
$ws.Range("D16").Value = '20.83'
$ws.Range("E16").Value = '  +4.59%  '
$ws.Range("D17").Value = '3.424.41'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("E18").Value = '  +10.26%  '
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").Value = '61.989.28'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '444.48'
$ws.Range("E21").Value = '  +42.04%  '
$ws.Range("D22").Value = '90.70'
$ws.Range("E22").Value = '  +7.02%  '
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '13.04'
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +3.14%  '
$ws.Range("D26").Value = '33.74'
$ws.Range("E26").Value = '  +13.89%  '
$ws.Range("D27").Value = '8.83'
$ws.Range("E27").Value = '  +8.30%  '
$ws.Range("D28").Value = '4.75'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = '7.58'
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = '2.72'
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("D31").Value = '11.95'
$ws.Range("E31").Value = '  +5.55%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '0.170'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("D34").Value = '42.48'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '0.0500'
$ws.Range("E36").Value = '  +3.72%  '
$ws.Range("D37").Value = '53.87'
$ws.Range("E37").Value = '  +4.12%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("E40").Value = '  +7.99%  '
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").Value = '141.09'
$ws.Range("E43").Value = '  +2.66%  '
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("E46").Value = '  +8.86%  '
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").Value = '22.12'
$ws.Range("E48").Value = '  +3.88%  '
$ws.Range("D49").Value = '3.758.61'
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").Value = '2.110.65'
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = '2.35'
$ws.Range("E51").Value = '  +1.88%  '

# Restore the original (default) cell style on column D so the
# number-format-as-text trick above doesn't leave a stray style behind.
$priceRange.Style = "Normal"
